$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 340, pushing existing rows 340:363 down to 341:364.
$ws.Rows("340:340").Insert()

# Populate the newly inserted row 340 with the new price record.
$ws.Range("A340").Value = 4
$ws.Range("B340").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C340").Value = "Los Lagos"
$ws.Range("D340").Value = 44931
$ws.Range("E340").Value = 10
$ws.Range("F340").Value = 100112043
$ws.Range("G340").Value = "Pepino ensalada"
$ws.Range("H340").Value = "Sin especificar"
$ws.Range("I340").Value = "Primera"
$ws.Range("J340").Value = 200
$ws.Range("K340").Value = 20000
$ws.Range("L340").Value = 20000
$ws.Range("M340").Value = 20000
$ws.Range("N340").Value = '$/caja 60 unidades'
$ws.Range("O340").Value = "Región de Arica y Parinacota"
$ws.Range("P340").Value = 333
$ws.Range("Q340").Value = 60
$ws.Range("R340").Value = "Hortaliza"
